$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2017")

# Row 2 (ADM): updated cost-basis price and June dividend
$ws.Range("D2").Value = 27.971
$ws.Range("M2").Value = 8.8800000000000008

# Row 3 (EMR): updated cost-basis price and June dividend
$ws.Range("D3").Value = 14.606999999999999
$ws.Range("M3").Value = 6.95

# Row 6 (KO): updated cost-basis price and June dividend
$ws.Range("D6").Value = 25.283000000000001
$ws.Range("M6").Value = 3.02

# Row 11 (UL): updated cost-basis price and June dividend
$ws.Range("D11").Value = 15.205
$ws.Range("M11").Value = 5.7

# Update the active selection shown in the sheet view
$ws.Activate()
$ws.Range("J18").Select()
